# Apply the "Improve video suggestion layout" edit:
#  - Replace the 4 existing sample rows (2-5) with new "economics"
#    themed video rows, and append 4 brand-new rows (6-9) with the
#    same shape, extending the sheet dimension from A1:S5 to A1:S9.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 6:9 need the same date display format as the existing
# "published" column (C2:C5) so the underlying style (s="2") matches.
$ws.Range("C6:C9").NumberFormat = $ws.Range("C2").NumberFormat

# Row 2
$ws.Range("A2").Value = 'How The Economic Machine Works by Ray Dalio'
$ws.Range("B2").Value = 'Principles by Ray Dalio'
$ws.Range("C2").Value = 41539.75540509259
$ws.Range("D2").Value = 43654314
$ws.Range("E2").Value = 671653
$ws.Range("F2").Value = 24368
$ws.Range("G2").Value = 31
$ws.Range("H2").Value = '31:00'
$ws.Range("I2").Value = 'PHe0bXAIuk0'
$ws.Range("J2").Value = 'https://www.youtube.com/watch?v=PHe0bXAIuk0'
$ws.Range("K2").Value = 0.01538571880891314
$ws.Range("L2").Value = 783.5369774919614
$ws.Range("M2").Value = 10271.60329411765
$ws.Range("N2").Value = 0.170419573613995
$ws.Range("O2").Value = 0.4489254980311678
$ws.Range("P2").Value = 0.9999999999999019
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 6.409109716904027
$ws.Range("S2").Value = 1

# Row 3
$ws.Range("A3").Value = 'Every Major Economic Theory Explained in 20 Minutes'
$ws.Range("B3").Value = 'Adam''s Axiom'
$ws.Range("C3").Value = 45649.54217592593
$ws.Range("D3").Value = 498390
$ws.Range("E3").Value = 18280
$ws.Range("F3").Value = 507
$ws.Range("G3").Value = 20.62
$ws.Range("H3").Value = '20:37'
$ws.Range("I3").Value = 'dQ_UPQa3CUE'
$ws.Range("J3").Value = 'https://www.youtube.com/watch?v=dQ_UPQa3CUE'
$ws.Range("K3").Value = 0.03667810349324826
$ws.Range("L3").Value = 24.46911196911197
$ws.Range("M3").Value = 3534.68085106383
$ws.Range("N3").Value = 0.5211902786773303
$ws.Range("O3").Value = 0.01401951483167455
$ws.Range("P3").Value = 0.3392962401352193
$ws.Range("Q3").Value = 0.01138080102074723
$ws.Range("R3").Value = 2.632260188142492
$ws.Range("S3").Value = 2

# Row 4
$ws.Range("A4").Value = 'Intro to Economics: Crash Course Econ #1'
$ws.Range("B4").Value = 'CrashCourse'
$ws.Range("C4").Value = 42193.94049768519
$ws.Range("D4").Value = 8264493
$ws.Range("E4").Value = 116672
$ws.Range("F4").Value = 4446
$ws.Range("G4").Value = 12.15
$ws.Range("H4").Value = '12:09'
$ws.Range("I4").Value = '3ez10ADR_gM'
$ws.Range("J4").Value = 'https://www.youtube.com/watch?v=3ez10ADR_gM'
$ws.Range("K4").Value = 0.01411726042964765
$ws.Range("L4").Value = 362.9387755102041
$ws.Range("M4").Value = 2298.246106785317
$ws.Range("N4").Value = 0.1495229926703963
$ws.Range("O4").Value = 0.2079448388923191
$ws.Range("P4").Value = 0.2180365564625309
$ws.Range("Q4").Value = 0.1892872809526882
$ws.Range("R4").Value = 1.897142887088796
$ws.Range("S4").Value = 3

# Row 5
$ws.Range("A5").Value = 'What Everyone Gets Wrong About Global Debt | Economics Explained'
$ws.Range("B5").Value = 'Economics Explained'
$ws.Range("C5").Value = 45109.56671296297
$ws.Range("D5").Value = 1641344
$ws.Range("E5").Value = 29325
$ws.Range("F5").Value = 1392
$ws.Range("G5").Value = 16.68
$ws.Range("H5").Value = '16:41'
$ws.Range("I5").Value = 'IAqj30s4lH8'
$ws.Range("J5").Value = 'https://www.youtube.com/watch?v=IAqj30s4lH8'
$ws.Range("K5").Value = 0.01786645578257818
$ws.Range("L5").Value = 82.95589988081048
$ws.Range("M5").Value = 2410.196769456682
$ws.Range("N5").Value = 0.2112872298314481
$ws.Range("O5").Value = 0.04752936968950978
$ws.Range("P5").Value = 0.2290157868710058
$ws.Range("Q5").Value = 0.03756367843869182
$ws.Range("R5").Value = 1.491095146363765
$ws.Range("S5").Value = 4

# Row 6
$ws.Range("A6").Value = 'The Most Important Economic Schools of Thought | Economics Explained'
$ws.Range("B6").Value = 'Economics Explained'
$ws.Range("C6").Value = 44084.53758101852
$ws.Range("D6").Value = 1841614
$ws.Range("E6").Value = 45281
$ws.Range("F6").Value = 3174
$ws.Range("G6").Value = 26.08
$ws.Range("H6").Value = '26:05'
$ws.Range("I6").Value = 'o6UXRZ2XwgU'
$ws.Range("J6").Value = 'https://www.youtube.com/watch?v=o6UXRZ2XwgU'
$ws.Range("K6").Value = 0.02458767146644194
$ws.Range("L6").Value = 121.2375859434683
$ws.Range("M6").Value = 1079.492379835873
$ws.Range("N6").Value = 0.3220125243094374
$ws.Range("O6").Value = 0.0694627633580017
$ws.Range("P6").Value = 0.09851088631401388
$ws.Range("Q6").Value = 0.04215147887553508
$ws.Range("R6").Value = 1.484798716337427
$ws.Range("S6").Value = 5

# Row 7
$ws.Range("A7").Value = 'How The Economy Works For DUMMIES: Global Economics 101 -Robert Kiyosaki'
$ws.Range("B7").Value = 'The Rich Dad Channel'
$ws.Range("C7").Value = 43585.81953703704
$ws.Range("D7").Value = 165339
$ws.Range("E7").Value = 3975
$ws.Range("F7").Value = 215
$ws.Range("G7").Value = 3.5
$ws.Range("H7").Value = '3:30'
$ws.Range("I7").Value = '9iV55N2nuJY'
$ws.Range("J7").Value = 'https://www.youtube.com/watch?v=9iV55N2nuJY'
$ws.Range("K7").Value = 0.02404151470614918
$ws.Range("L7").Value = 59.72222222222222
$ws.Range("M7").Value = 75.01769509981851
$ws.Range("N7").Value = 0.3130151388721953
$ws.Range("O7").Value = 0.03421769377172024
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0.003751243307205069
$ws.Range("R7").Value = 1.014983290774436
$ws.Range("S7").Value = 6

# Row 8
$ws.Range("A8").Value = 'Introduction to Economics Part 1 - Professor Ryan'
$ws.Range("B8").Value = 'Prof Ryan'
$ws.Range("C8").Value = 43613.19787037037
$ws.Range("D8").Value = 189733
$ws.Range("E8").Value = 4337
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 17.73
$ws.Range("H8").Value = '17:44'
$ws.Range("I8").Value = 'qnEZQRpWWi8'
$ws.Range("J8").Value = 'https://www.youtube.com/watch?v=qnEZQRpWWi8'
$ws.Range("K8").Value = 0.02285843791011579
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 87.15342214056041
$ws.Range("N8").Value = 0.2935251339694559
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0.001190175566408197
$ws.Range("Q8").Value = 0.004310062920009556
$ws.Range("R8").Value = 0.8927660544476114
$ws.Range("S8").Value = 7

# Row 9
$ws.Range("A9").Value = 'Thomas Sowell -- Basic Economics'
$ws.Range("B9").Value = 'Hoover Institution'
$ws.Range("C9").Value = 40549.97076388889
$ws.Range("D9").Value = 2520647
$ws.Range("E9").Value = 44076
$ws.Range("F9").Value = 1256
$ws.Range("G9").Value = 33.53
$ws.Range("H9").Value = '33:32'
$ws.Range("I9").Value = 'bOMksnSaAJ4'
$ws.Range("J9").Value = 'https://www.youtube.com/watch?v=bOMksnSaAJ4'
$ws.Range("K9").Value = 0.01748598673277139
$ws.Range("L9").Value = 37.34760630389533
$ws.Range("M9").Value = 481.0395038167939
$ws.Range("N9").Value = 0.2050193834237279
$ws.Range("O9").Value = 0.02139821507073698
$ws.Range("P9").Value = 0.03981938902724899
$ws.Range("Q9").Value = 0.05770681863701207
$ws.Range("R9").Value = 0.8927263847684287
$ws.Range("S9").Value = 8
